# Update gh-pages to output generated at 456a3b4
# Applies the numeric "想去人数" (want-to-go count) bumps, a couple of
# "最低票价" (min price) cells flipping from a number to a sold-out /
# unavailable label, and refreshes the "全部类型" (all types) sheet's
# rolling window of 展览+演出 events for 2024-08-10/11 (new event added,
# oldest one dropped, remaining ones shift up one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 31
$ws1.Range("F6").Value  = 711
$ws1.Range("F9").Value  = 877
$ws1.Range("F11").Value = 273
$ws1.Range("F13").Value = 387
$ws1.Range("F14").Value = 732
$ws1.Range("F15").Value = 1053
$ws1.Range("F16").Value = 12209
$ws1.Range("F17").Value = 663
$ws1.Range("F18").Value = 56
$ws1.Range("F20").Value = 47
$ws1.Range("F23").Value = 1808
$ws1.Range("F27").Value = 195
$ws1.Range("F28").Value = 112
$ws1.Range("F29").Value = 312
$ws1.Range("F31").Value = 277
$ws1.Range("F32").Value = 92
$ws1.Range("F37").Value = 1214

# ---------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value  = "不可售"
$ws2.Range("F5").Value  = 107
$ws2.Range("F6").Value  = 160
$ws2.Range("F9").Value  = 255
$ws2.Range("F10").Value = 4442
$ws2.Range("F14").Value = 67

# ---------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 842

# ---------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 842
$ws4.Range("F3").Value = 31

# Row 4: was 广州·七夕《梁祝》中外经典名曲音乐会 -> becomes 广州·火影忍者only
$ws4.Range("B4").Value = "2024-08-10"
$ws4.Range("C4").Value = "广州·火影忍者only"
$ws4.Range("D4").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws4.Range("E4").Value = "2024.08.10 10:00-08.10 17:00"
$ws4.Range("F4").Value = 1267
$ws4.Range("G4").Value = 70
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=85704"
$ws4.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202408/49fmnEM71723086988340.jpeg"

# Row 5: becomes 广州·电子音乐仓库派对：RoaringLand咆哮大陆
$ws4.Range("B5").Value = "2024-08-10"
$ws4.Range("C5").Value = "广州·电子音乐仓库派对：RoaringLand咆哮大陆"
$ws4.Range("D5").Value = "革新路124号太古仓码头4号仓 MAO Livehouse 广州(太古仓店)"
$ws4.Range("E5").Value = "2024.08.10 23:00-08.11 04:00"
$ws4.Range("F5").Value = 4
$ws4.Range("G5").Value = 230
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=90061"
$ws4.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202407/HRdYR5SK1721892863746.png"

# Row 6: becomes 广州·系统任务：重生之我是音乐一体机！王子健2024巡回演出
$ws4.Range("B6").Value = "2024-08-10"
$ws4.Range("C6").Value = "广州·系统任务：重生之我是音乐一体机！王子健2024巡回演出"
$ws4.Range("D6").Value = "南洲路154号 SD Livehouse"
$ws4.Range("E6").Value = "2024.08.10 20:00-08.10 22:00"
$ws4.Range("F6").Value = 107
$ws4.Range("G6").Value = 128
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=87585"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202407/TxCZBf4D1721810695745.png"

# Row 7: becomes 广州·OVO动漫嘉年华2th (also date bumps to 2024-08-11)
$ws4.Range("B7").Value = "2024-08-11"
$ws4.Range("C7").Value = "广州·OVO动漫嘉年华2th"
$ws4.Range("D7").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws4.Range("E7").Value = "2024.08.11 10:00-08.11 17:00"
$ws4.Range("F7").Value = 877
$ws4.Range("G7").Value = 45
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=89822"
$ws4.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202407/dMzUkTL41721797711362.jpeg"

# Row 8: becomes 广州·咒术回战ONLY
$ws4.Range("B8").Value = "2024-08-11"
$ws4.Range("C8").Value = "广州·咒术回战ONLY"
$ws4.Range("D8").Value = "西环路1号 广州岭南会展中心"
$ws4.Range("E8").Value = "2024.08.11 10:00-08.11 17:00"
$ws4.Range("F8").Value = 729
$ws4.Range("G8").Value = 60
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=87433"
$ws4.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202406/kNv9yqGn1718350051848.jpeg"

# Row 9: becomes 广州·妖都偶像梦幻祭only5.0
$ws4.Range("B9").Value = "2024-08-11"
$ws4.Range("C9").Value = "广州·妖都偶像梦幻祭only5.0"
$ws4.Range("D9").Value = "同泰路颐和山庄内国际会议厅 颐和山庄"
$ws4.Range("E9").Value = "2024.08.11 09:30-08.11 18:00"
$ws4.Range("F9").Value = 273
$ws4.Range("G9").Value = 68
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=89150"
$ws4.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202407/HSCluzha1719829266544.jpeg"

# Row 10: becomes 广州·虚拟主播Virtual Only
$ws4.Range("B10").Value = "2024-08-11"
$ws4.Range("C10").Value = "广州·虚拟主播Virtual Only"
$ws4.Range("D10").Value = "逸景路462号珠江国际纺织城d区6层 珠江时尚馆"
$ws4.Range("E10").Value = "2024.08.11 10:00-08.11 17:30"
$ws4.Range("F10").Value = 387
$ws4.Range("G10").Value = 80
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=88934"
$ws4.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202407/XTc8Vx4z1720443835316.jpeg"

# Row 11: new event 广州·迷宫饭only (replaces 广州·虚拟主播Virtual Only, which moved to row 10)
$ws4.Range("B11").Value = "2024-08-11"
$ws4.Range("C11").Value = "广州·迷宫饭only"
$ws4.Range("D11").Value = "大石街石北工业大道644号 巨大创意产业园"
$ws4.Range("E11").Value = "2024.08.11 11:00-08.11 17:00"
$ws4.Range("F11").Value = 732
$ws4.Range("G11").Value = "已售罄"
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=88675"
$ws4.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202407/TMAyXTkr1720077147308.png"

# Remaining independent "想去人数" bumps further down the same sheet
$ws4.Range("F14").Value = 1053
$ws4.Range("F15").Value = 12210
$ws4.Range("F16").Value = 255
$ws4.Range("F17").Value = 663
$ws4.Range("F18").Value = 56
$ws4.Range("F20").Value = 47
$ws4.Range("F22").Value = 1808
$ws4.Range("F25").Value = 195
$ws4.Range("F29").Value = 112
$ws4.Range("F31").Value = 67
$ws4.Range("F33").Value = 312
$ws4.Range("F36").Value = 277
$ws4.Range("F37").Value = 93
$ws4.Range("F45").Value = 1214
